# Apply the "Atualizacoes 16 de janeiro de 2024" edits.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$incl = $wb.Worksheets.Item("Include from ")

# Rename the "Include from " sheet to reflect the new Title value.
$incl.Name = "Include from Duration of Stro"

# Metadata sheet updates
$meta.Range("B3").Value = "0.0.0"
$meta.Range("B5").Value = "Duration of Stroke Symptoms"
$meta.Range("B7").Value = "'false"
$meta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Include sheet updates
$incl.Range("B6").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/SymtDurCS"
